$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "FORMATO X" label that used to live in J4 (the row had no
# other content, so clearing it drops the row entirely)
$ws.Range("J4").Clear()

# Insert a new blank row above the "Fecha:" row, pushing the Fecha/Hora
# labels and the table header down by one row
$ws.Rows("7:7").Insert()

# The inserted row inherited the formatting of the row above it (row 6);
# trim it back so only the J column keeps the small right-hand box style
# that matches the Fecha/Hora column below it
$ws.Range("A7:I7").Clear()
$ws.Range("J8").Copy()
$ws.Range("J7").PasteSpecial(-4122)
$ws.Rows("7:7").RowHeight = 15.75
$excel.CutCopyMode = $false

# Rename the last header column from "No. de Proyecto" to "Clave de Proyecto"
$ws.Range("E11").Value = "Clave de Proyecto"

# Match the saved selection state
$ws.Range("J4").Select()
